# "added representants for types"
#
# The diff widens column A (so the longer Wuchsgebiet / type labels are
# fully visible) and leaves the selection parked on B19 - the row for the
# representative "13: Ostniedersächsisches Tiefland" region - after the
# author picked representative rows/types. No cell values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the (now longer) type/region labels.
# NOTE: the target OOXML width is 56.7109375 character-units; this host
# quantizes ColumnWidth to whole "pixels" (1/6ths) internally, so
# 56.7109375 itself is not exactly reachable - 55.8333333333 is the
# closest input that round-trips to the nearest reachable width.
$ws.Range("A1").EntireColumn.ColumnWidth = 55.8333333333

# Leave the cursor on B19, matching the saved selection in the sheet view.
$ws.Range("B19").Select()
